# Rename the sheet "Property1" to "DataNode" as part of unifying the
# DataNode / DataTable / Entity concept (see commit message), and restore
# the last active selection to E23 on that sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "DataNode"
$ws.Range("E23").Select()
